# Updated symbol list on Tue Jan 31 08:16:03 UTC 2023 with GitHub Actions
# Applies the refreshed cryptocurrency price/volume/hour data onto Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the cell's exact string
# representation (e.g. leading zeros, trailing zeros, "%" suffixes) instead
# of letting Excel auto-convert numeric-looking strings into numbers/percents.
# NumberFormat is temporarily switched to Text ("@") for the write, then the
# cell's original style object is restored so no formatting is left behind.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$updates = @(
    @{ Row = 2; Col = "D"; Value = "311.74" },
    @{ Row = 2; Col = "E"; Value = "-0.55%" },
    @{ Row = 2; Col = "G"; Value = "8" },
    @{ Row = 3; Col = "D"; Value = "37.66" },
    @{ Row = 3; Col = "E"; Value = "-4.13%" },
    @{ Row = 3; Col = "G"; Value = "8" },
    @{ Row = 4; Col = "D"; Value = "5.065" },
    @{ Row = 4; Col = "E"; Value = "-1.36%" },
    @{ Row = 4; Col = "G"; Value = "8" },
    @{ Row = 5; Col = "D"; Value = "0.07783" },
    @{ Row = 5; Col = "E"; Value = "-4.00%" },
    @{ Row = 5; Col = "G"; Value = "8" },
    @{ Row = 6; Col = "D"; Value = "4.358" },
    @{ Row = 6; Col = "E"; Value = "-2.67%" },
    @{ Row = 6; Col = "G"; Value = "8" },
    @{ Row = 7; Col = "B"; Value = "KuCoinToken" },
    @{ Row = 7; Col = "C"; Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs" },
    @{ Row = 7; Col = "D"; Value = "8.220" },
    @{ Row = 7; Col = "E"; Value = "-0.96%" },
    @{ Row = 7; Col = "G"; Value = "8" },
    @{ Row = 8; Col = "B"; Value = "FTXToken" },
    @{ Row = 8; Col = "C"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" },
    @{ Row = 8; Col = "D"; Value = "1.888" },
    @{ Row = 8; Col = "E"; Value = "-3.70%" },
    @{ Row = 8; Col = "G"; Value = "8" },
    @{ Row = 9; Col = "E"; Value = "-10.10%" },
    @{ Row = 9; Col = "G"; Value = "8" },
    @{ Row = 10; Col = "D"; Value = "0.9181" },
    @{ Row = 10; Col = "E"; Value = "-2.12%" },
    @{ Row = 10; Col = "G"; Value = "8" },
    @{ Row = 11; Col = "D"; Value = "0.1199" },
    @{ Row = 11; Col = "E"; Value = "-9.81%" },
    @{ Row = 11; Col = "G"; Value = "8" },
    @{ Row = 12; Col = "D"; Value = "0.1927" },
    @{ Row = 12; Col = "E"; Value = "-1.97%" },
    @{ Row = 12; Col = "G"; Value = "8" },
    @{ Row = 13; Col = "D"; Value = "0.09106" },
    @{ Row = 13; Col = "E"; Value = "0.37%" },
    @{ Row = 13; Col = "G"; Value = "8" },
    @{ Row = 14; Col = "D"; Value = "0.03405" },
    @{ Row = 14; Col = "E"; Value = "-2.24%" },
    @{ Row = 14; Col = "G"; Value = "8" },
    @{ Row = 15; Col = "D"; Value = "0.09705" },
    @{ Row = 15; Col = "E"; Value = "-0.04%" },
    @{ Row = 15; Col = "G"; Value = "8" },
    @{ Row = 16; Col = "D"; Value = "0.001363" },
    @{ Row = 16; Col = "E"; Value = "-2.97%" },
    @{ Row = 16; Col = "G"; Value = "8" },
    @{ Row = 17; Col = "D"; Value = "0.005868" },
    @{ Row = 17; Col = "E"; Value = "-1.75%" },
    @{ Row = 17; Col = "G"; Value = "8" },
    @{ Row = 18; Col = "E"; Value = "0.02%" },
    @{ Row = 18; Col = "G"; Value = "8" },
    @{ Row = 19; Col = "D"; Value = "0.3411" },
    @{ Row = 19; Col = "E"; Value = "-1.36%" },
    @{ Row = 19; Col = "G"; Value = "8" },
    @{ Row = 20; Col = "D"; Value = "5.133" },
    @{ Row = 20; Col = "E"; Value = "2.52%" },
    @{ Row = 20; Col = "G"; Value = "8" },
    @{ Row = 21; Col = "D"; Value = "0.1268" },
    @{ Row = 21; Col = "E"; Value = "-1.73%" },
    @{ Row = 21; Col = "G"; Value = "8" },
    @{ Row = 22; Col = "E"; Value = "3.80%" },
    @{ Row = 22; Col = "G"; Value = "8" },
    @{ Row = 23; Col = "D"; Value = "0.02103" },
    @{ Row = 23; Col = "E"; Value = "5,585.46%" },
    @{ Row = 23; Col = "G"; Value = "8" },
    @{ Row = 24; Col = "D"; Value = "0.04374" },
    @{ Row = 24; Col = "E"; Value = "0.09%" },
    @{ Row = 24; Col = "G"; Value = "8" },
    @{ Row = 25; Col = "D"; Value = "0.001212" },
    @{ Row = 25; Col = "E"; Value = "-2.70%" },
    @{ Row = 25; Col = "G"; Value = "8" },
    @{ Row = 26; Col = "D"; Value = "0.004269" },
    @{ Row = 26; Col = "E"; Value = "-9.83%" },
    @{ Row = 26; Col = "G"; Value = "8" },
    @{ Row = 27; Col = "E"; Value = "-66.61%" },
    @{ Row = 27; Col = "G"; Value = "8" },
    @{ Row = 28; Col = "G"; Value = "8" },
    @{ Row = 29; Col = "G"; Value = "8" },
    @{ Row = 30; Col = "G"; Value = "8" },
    @{ Row = 31; Col = "G"; Value = "8" },
    @{ Row = 32; Col = "G"; Value = "8" },
    @{ Row = 33; Col = "G"; Value = "8" },
    @{ Row = 34; Col = "G"; Value = "8" },
    @{ Row = 35; Col = "G"; Value = "8" },
    @{ Row = 36; Col = "G"; Value = "8" },
    @{ Row = 37; Col = "G"; Value = "8" },
    @{ Row = 38; Col = "G"; Value = "8" },
    @{ Row = 39; Col = "D"; Value = "0.02119" },
    @{ Row = 39; Col = "E"; Value = "-4.21%" },
    @{ Row = 39; Col = "G"; Value = "8" },
    @{ Row = 40; Col = "D"; Value = "0.04973" },
    @{ Row = 40; Col = "E"; Value = "-4.59%" },
    @{ Row = 40; Col = "G"; Value = "8" },
    @{ Row = 41; Col = "D"; Value = "0.007799" },
    @{ Row = 41; Col = "E"; Value = "2.37%" },
    @{ Row = 41; Col = "G"; Value = "8" },
    @{ Row = 42; Col = "D"; Value = "0.009924" },
    @{ Row = 42; Col = "E"; Value = "-3.94%" },
    @{ Row = 42; Col = "G"; Value = "8" },
    @{ Row = 43; Col = "E"; Value = "-3.34%" },
    @{ Row = 43; Col = "G"; Value = "8" },
    @{ Row = 44; Col = "D"; Value = "0.002059" },
    @{ Row = 44; Col = "E"; Value = "-1.98%" },
    @{ Row = 44; Col = "G"; Value = "8" },
    @{ Row = 45; Col = "D"; Value = "0.008808" },
    @{ Row = 45; Col = "E"; Value = "-3.19%" },
    @{ Row = 45; Col = "G"; Value = "8" },
    @{ Row = 46; Col = "D"; Value = "0.00006665" },
    @{ Row = 46; Col = "E"; Value = "0.93%" },
    @{ Row = 46; Col = "G"; Value = "8" },
    @{ Row = 47; Col = "D"; Value = "0.00000000750" },
    @{ Row = 47; Col = "E"; Value = "-0.10%" },
    @{ Row = 47; Col = "G"; Value = "8" },
    @{ Row = 48; Col = "D"; Value = "0.002915" },
    @{ Row = 48; Col = "E"; Value = "-3.24%" },
    @{ Row = 48; Col = "G"; Value = "8" },
    @{ Row = 49; Col = "D"; Value = "0.001200" },
    @{ Row = 49; Col = "E"; Value = "-29.06%" },
    @{ Row = 49; Col = "G"; Value = "8" },
    @{ Row = 50; Col = "D"; Value = "0.00002099" },
    @{ Row = 50; Col = "E"; Value = "-0.10%" },
    @{ Row = 50; Col = "G"; Value = "8" },
    @{ Row = 51; Col = "D"; Value = "0.0001999" },
    @{ Row = 51; Col = "E"; Value = "-0.10%" },
    @{ Row = 51; Col = "G"; Value = "8" }
)

foreach ($u in $updates) {
    $addr = "$($u.Col)$($u.Row)"
    $range = $ws.Range($addr)
    Set-TextValue $range $u.Value
}
